$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1
$ws = $excel.ActiveWorkbook.ActiveSheet
$ws.Range("C49").Select() | Out-Null
